$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

function Set-TextId($rangeAddr, $text) {
    # Writing a numeric-looking string through .Value coerces it to a
    # number, so build it as a text formula first and then flatten the
    # formula down to a plain (text) value via copy/paste-special. This
    # keeps the cell as a shared string (same as the original file) and
    # does not touch the cell's style.
    $r = $ws.Range($rangeAddr)
    $r.Formula = '="' + $text + '"'
    $r.Copy() | Out-Null
    $r.PasteSpecial($xlPasteValues) | Out-Null
}

# id column (B) - reassign each fixture's id to its new row position
Set-TextId "B194" "8209691"
Set-TextId "B195" "8209690"
Set-TextId "B196" "8209692"
Set-TextId "B197" "8209693"
Set-TextId "B198" "8209694"
Set-TextId "B199" "8209689"

# Row 194 - Mezokovesd Zsory vs MTK Budapest
$ws.Range("D194").Value = 45429.64583333334
$ws.Range("E194").Value = "Mezokovesd Zsory"
$ws.Range("F194").Value = "MTK Budapest"
$ws.Range("J194").Value = 3
$ws.Range("K194").Value = 3.6
$ws.Range("L194").Value = 2.1
$ws.Range("M194").Value = 3.1
$ws.Range("N194").Value = 3.6
$ws.Range("O194").Value = 2.05
$ws.Range("P194").Value = 0.25
$ws.Range("Q194").Value = 2
$ws.Range("R194").Value = 1.85
$ws.Range("S194").Value = 3.25
$ws.Range("T194").Value = 1.975
$ws.Range("U194").Value = 1.875

# Row 195 - Paksi vs Kisvarda FC
$ws.Range("D195").Value = 45430.41666666666
$ws.Range("E195").Value = "Paksi"
$ws.Range("F195").Value = "Kisvarda FC"
$ws.Range("J195").Value = 1.444
$ws.Range("K195").Value = 4.333
$ws.Range("L195").Value = 6
$ws.Range("M195").Value = 1.42
$ws.Range("N195").Value = 4.75
$ws.Range("O195").Value = 5.5
$ws.Range("P195").Value = -1.25
$ws.Range("Q195").Value = 1.975
$ws.Range("R195").Value = 1.875
$ws.Range("S195").Value = 3
$ws.Range("T195").Value = 1.825
$ws.Range("U195").Value = 2.025

# Row 196 - Puskas Academy vs Debreceni VSC
$ws.Range("D196").Value = 45430.41666666666
$ws.Range("E196").Value = "Puskas Academy"
$ws.Range("F196").Value = "Debreceni VSC"
$ws.Range("J196").Value = 1.7
$ws.Range("K196").Value = 3.6
$ws.Range("L196").Value = 4.5
$ws.Range("M196").Value = 1.5
$ws.Range("N196").Value = 4
$ws.Range("O196").Value = 5.75
$ws.Range("P196").Value = -1
$ws.Range("Q196").Value = 1.85
$ws.Range("R196").Value = 2
$ws.Range("S196").Value = 3
$ws.Range("T196").Value = 2.025
$ws.Range("U196").Value = 1.825

# Row 197 - MOL Fehervar FC vs Diosgyori VTK
$ws.Range("D197").Value = 45430.41666666666
$ws.Range("E197").Value = "MOL Fehervar FC"
$ws.Range("F197").Value = "Diosgyori VTK"
$ws.Range("J197").Value = 1.571
$ws.Range("K197").Value = 4
$ws.Range("L197").Value = 5
$ws.Range("M197").Value = 1.75
$ws.Range("N197").Value = 3.75
$ws.Range("O197").Value = 4
$ws.Range("P197").Value = -0.75
$ws.Range("Q197").Value = 2.025
$ws.Range("R197").Value = 1.825
$ws.Range("S197").Value = 3
$ws.Range("T197").Value = 1.875
$ws.Range("U197").Value = 1.975

# Row 198 - Kecskemeti TE vs Zalaegerszegi TE
$ws.Range("D198").Value = 45431.51041666666
$ws.Range("E198").Value = "Kecskemeti TE"
$ws.Range("F198").Value = "Zalaegerszegi TE"
$ws.Range("J198").Value = 2.05
$ws.Range("K198").Value = 3.5
$ws.Range("L198").Value = 3.2
$ws.Range("M198").Value = 1.95
$ws.Range("N198").Value = 3.7
$ws.Range("O198").Value = 3.4
$ws.Range("P198").Value = -0.5
$ws.Range("Q198").Value = 2.025
$ws.Range("R198").Value = 1.825
$ws.Range("S198").Value = 2.75
$ws.Range("T198").Value = 2
$ws.Range("U198").Value = 1.85

# Row 199 - Ferencvarosi TC vs Ujpest
$ws.Range("D199").Value = 45431.63541666666
$ws.Range("E199").Value = "Ferencvarosi TC"
$ws.Range("F199").Value = "Ujpest"
$ws.Range("J199").Value = 1.222
$ws.Range("K199").Value = 5.25
$ws.Range("L199").Value = 12
$ws.Range("M199").Value = 1.3
$ws.Range("N199").Value = 4.75
$ws.Range("O199").Value = 9
$ws.Range("P199").Value = -1.5
$ws.Range("Q199").Value = 1.875
$ws.Range("R199").Value = 1.975
$ws.Range("S199").Value = 3
$ws.Range("T199").Value = 1.925
$ws.Range("U199").Value = 1.925
